$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 612 ("「これ、私の眼鏡じゃありません」...") and shift subsequent rows up.
$ws.Rows.Item(612).Delete()
